$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 1.53
$ws.Range("J2").Value = 7.5
$ws.Range("L2").Value = 2.2
$ws.Range("N2").Value = 7.5
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.5
$ws.Range("X2").Value = 34
$ws.Range("Z2").Value = 81
$ws.Range("AD2").Value = 8
$ws.Range("AE2").Value = 26
$ws.Range("AF2").Value = 101
$ws.Range("AK2").Value = 10
$ws.Range("AN2").Value = 8
$ws.Range("AP2").Value = 51
$ws.Range("AQ2").Value = 201
$ws.Range("AR2").Value = 251
$ws.Range("AW2").Value = 3.2
$ws.Range("AX2").Value = 8
$ws.Range("AZ2").Value = 26

# Row 7
$ws.Range("S7").Value = 1.41
$ws.Range("T7").Value = 2.62

# Row 9
$ws.Range("S9").Value = 1.37

# Row 10
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 6
$ws.Range("O10").Value = 1.53
$ws.Range("P10").Value = 2.38
$ws.Range("R10").Value = 1.41
